$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.056748553031213
$ws.Range("D2").Value = 1.054271244083026
$ws.Range("E2").Value = 1.062098227298203
$ws.Range("F2").Value = 1.07078119812752
$ws.Range("I2").Value = 1.041198566114515
$ws.Range("J2").Value = 1.061748569773783
$ws.Range("K2").Value = 1.057014914388477
$ws.Range("L2").Value = 1.064820496232767
$ws.Range("M2").Value = 1.073480121093358
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058371590284792
$ws.Range("D3").Value = 1.055503062703237
$ws.Range("E3").Value = 1.063544809407953
$ws.Range("F3").Value = 1.072313733926125
$ws.Range("I3").Value = 1.041579471626975
$ws.Range("J3").Value = 1.063020560694668
$ws.Range("K3").Value = 1.058058873017664
$ws.Range("L3").Value = 1.066080238859819
$ws.Range("M3").Value = 1.074827314628469
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059419962158906
$ws.Range("D4").Value = 1.056298266789983
$ws.Range("E4").Value = 1.064479359051927
$ws.Range("F4").Value = 1.073303973633464
$ws.Range("I4").Value = 1.041823749525427
$ws.Range("J4").Value = 1.063841422877484
$ws.Range("K4").Value = 1.058731960273041
$ws.Range("L4").Value = 1.066893373492229
$ws.Range("M4").Value = 1.075697122019768
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059860266817201
$ws.Range("D5").Value = 1.056632130791778
$ws.Range("E5").Value = 1.064871896053509
$ws.Range("F5").Value = 1.073719940487428
$ws.Range("I5").Value = 1.041925921480718
$ws.Range("J5").Value = 1.064185994068419
$ws.Range("K5").Value = 1.059014352150664
$ws.Range("L5").Value = 1.067234743140965
$ws.Range("M5").Value = 1.076062337903095
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.059934170930517
$ws.Range("D6").Value = 1.056688162429831
$ws.Range("E6").Value = 1.064937784587706
$ws.Range("F6").Value = 1.073789764050708
$ws.Range("I6").Value = 1.041943046039831
$ws.Range("J6").Value = 1.064243818858655
$ws.Range("K6").Value = 1.059061733510471
$ws.Range("L6").Value = 1.067292033133768
$ws.Range("M6").Value = 1.076123633026365
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059425847210402
$ws.Range("D7").Value = 1.056302729615854
$ws.Range("E7").Value = 1.064484605505458
$ws.Range("F7").Value = 1.073309533087467
$ws.Range("I7").Value = 1.041825116801328
$ws.Range("J7").Value = 1.063846029084298
$ws.Range("K7").Value = 1.058735735853815
$ws.Range("L7").Value = 1.066897936735087
$ws.Range("M7").Value = 1.075702003814111
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05729745348192
$ws.Range("D8").Value = 1.054687933129256
$ws.Range("E8").Value = 1.062587418960484
$ws.Range("F8").Value = 1.071299423544586
$ws.Range("I8").Value = 1.041327750677337
$ws.Range("J8").Value = 1.062178905413878
$ws.Range("K8").Value = 1.057368230733854
$ws.Range("L8").Value = 1.065246651367257
$ws.Range("M8").Value = 1.073935813056144
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.053532394482879
$ws.Range("D9").Value = 1.051827871752655
$ws.Range("E9").Value = 1.059232594492014
$ws.Range("F9").Value = 1.067746147879225
$ws.Range("I9").Value = 1.040434409813155
$ws.Range("J9").Value = 1.059224019044551
$ws.Range("K9").Value = 1.054939666980678
$ws.Range("L9").Value = 1.062321201661766
$ws.Range("M9").Value = 1.070808528432475
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051011933803927
$ws.Range("D10").Value = 1.04991093595815
$ws.Range("E10").Value = 1.056987645181503
$ws.Range("F10").Value = 1.065369228472954
$ws.Range("I10").Value = 1.039827306375656
$ws.Range("J10").Value = 1.057242045399175
$ws.Range("K10").Value = 1.05330755772833
$ws.Range("L10").Value = 1.060359895997203
$ws.Range("M10").Value = 1.068713074750962
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0499179329975
$ws.Range("D11").Value = 1.049078360654492
$ws.Range("E11").Value = 1.056013451549505
$ws.Range("F11").Value = 1.064337962524582
$ws.Range("I11").Value = 1.039561648224686
$ws.Range("J11").Value = 1.056380862518025
$ws.Range("K11").Value = 1.052597648181521
$ws.Range("L11").Value = 1.059507913780696
$ws.Range("M11").Value = 1.06780309656176
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049511165769608
$ws.Range("D12").Value = 1.048768716818579
$ws.Range("E12").Value = 1.055651264925381
$ws.Range("F12").Value = 1.063954586920674
$ws.Range("I12").Value = 1.039462550214835
$ws.Range("J12").Value = 1.056060524425826
$ws.Range("K12").Value = 1.052333468324821
$ws.Range("L12").Value = 1.059191031299628
$ws.Range("M12").Value = 1.067464685007012
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049598437319534
$ws.Range("D13").Value = 1.048835154177492
$ws.Range("E13").Value = 1.055728970175658
$ws.Range("F13").Value = 1.064036836825271
$ws.Range("I13").Value = 1.039483826167673
$ws.Range("J13").Value = 1.05612925884063
$ws.Range("K13").Value = 1.052390158002105
$ws.Range("L13").Value = 1.059259022745975
$ws.Range("M13").Value = 1.067537293922256
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049884317857895
$ws.Range("D14").Value = 1.049052773367731
$ws.Range("E14").Value = 1.055983519833203
$ws.Range("F14").Value = 1.064306279115252
$ws.Range("I14").Value = 1.039553465354696
$ws.Range("J14").Value = 1.056354392637223
$ws.Range("K14").Value = 1.052575820989134
$ws.Range("L14").Value = 1.059481728763211
$ws.Range("M14").Value = 1.067775131677463
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050060404042584
$ws.Range("D15").Value = 1.049186804031039
$ws.Range("E15").Value = 1.056140312553986
$ws.Range("F15").Value = 1.064472249016425
$ws.Range("I15").Value = 1.039596316503418
$ws.Range("J15").Value = 1.056493044235096
$ws.Range("K15").Value = 1.052690149205888
$ws.Range("L15").Value = 1.059618889581002
$ws.Range("M15").Value = 1.067921617395446
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051084482688779
$ws.Range("D16").Value = 1.049966137241743
$ws.Range("E16").Value = 1.05705225378782
$ws.Range("F16").Value = 1.065437626243749
$ws.Range("I16").Value = 1.039844878412156
$ws.Range("J16").Value = 1.057299135783165
$ws.Range("K16").Value = 1.053354604101428
$ws.Range("L16").Value = 1.060416381085135
$ws.Range("M16").Value = 1.068773410718817
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05172614882067
$ws.Range("D17").Value = 1.05045430935014
$ws.Range("E17").Value = 1.057623717044401
$ws.Range("F17").Value = 1.066042627279511
$ws.Range("I17").Value = 1.040000048546624
$ws.Range("J17").Value = 1.057803972532112
$ws.Range("K17").Value = 1.053770538129589
$ws.Range("L17").Value = 1.060915890757344
$ws.Range("M17").Value = 1.069307006528785
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052100169409516
$ws.Range("D18").Value = 1.050738808389614
$ws.Range("E18").Value = 1.057956838424082
$ws.Range("F18").Value = 1.066395317623894
$ws.Range("I18").Value = 1.040090288816417
$ws.Range("J18").Value = 1.058098148843667
$ws.Range("K18").Value = 1.054012837622917
$ws.Range("L18").Value = 1.061206984099436
$ws.Range("M18").Value = 1.069617990550155
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052227658321593
$ws.Range("D19").Value = 1.050835774158416
$ws.Range("E19").Value = 1.058070389990805
$ws.Range("F19").Value = 1.066515542896171
$ws.Range("I19").Value = 1.040121013085403
$ws.Range("J19").Value = 1.058198407065318
$ws.Range("K19").Value = 1.054095403515766
$ws.Range("L19").Value = 1.061306195207369
$ws.Range("M19").Value = 1.069723985374168
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051657330314201
$ws.Range("D20").Value = 1.050401958348365
$ws.Range("E20").Value = 1.057562425562437
$ws.Range("F20").Value = 1.065977736796727
$ws.Range("I20").Value = 1.039983427988203
$ws.Range("J20").Value = 1.057749838021015
$ws.Range("K20").Value = 1.053725944236868
$ws.Range("L20").Value = 1.060862325283221
$ws.Range("M20").Value = 1.069249783026548
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049800144488591
$ws.Range("D21").Value = 1.048988700725993
$ws.Range("E21").Value = 1.055908570430238
$ws.Range("F21").Value = 1.064226943925455
$ws.Range("I21").Value = 1.039532969991311
$ws.Range("J21").Value = 1.05628810902122
$ws.Range("K21").Value = 1.052521161396493
$ws.Range("L21").Value = 1.05941615901408
$ws.Range("M21").Value = 1.067705105647619
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04863009931779
$ws.Range("D22").Value = 1.04809787883951
$ws.Range("E22").Value = 1.054866824953455
$ws.Range("F22").Value = 1.06312430870275
$ws.Range("I22").Value = 1.039247313202295
$ws.Range("J22").Value = 1.055366416507999
$ws.Range("K22").Value = 1.051760840205649
$ws.Range("L22").Value = 1.058504471856485
$ws.Range("M22").Value = 1.06673155686041
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049250590110704
$ws.Range("D23").Value = 1.048570336502217
$ws.Range("E23").Value = 1.055419257460373
$ws.Range("F23").Value = 1.063709014617792
$ws.Range("I23").Value = 1.039398977291677
$ws.Range("J23").Value = 1.055855277080088
$ws.Range("K23").Value = 1.052164171373782
$ws.Range("L23").Value = 1.058988007396318
$ws.Range("M23").Value = 1.067247879350656
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051688427215025
$ws.Range("D24").Value = 1.05042561426518
$ws.Range("E24").Value = 1.05759012117373
$ws.Range("F24").Value = 1.066007058620117
$ws.Range("I24").Value = 1.039990938930814
$ws.Range("J24").Value = 1.057774299961526
$ws.Range("K24").Value = 1.053746095248629
$ws.Range("L24").Value = 1.060886530025574
$ws.Range("M24").Value = 1.069275640647131
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.054507542535223
$ws.Range("D25").Value = 1.052569037721632
$ws.Range("E25").Value = 1.060101340116733
$ws.Range("F25").Value = 1.068666137540187
$ws.Range("I25").Value = 1.040667381931663
$ws.Range("J25").Value = 1.059990016044026
$ws.Range("K25").Value = 1.055569782746297
$ws.Range("L25").Value = 1.063079405866313
$ws.Range("M25").Value = 1.07161883803041
